# Apply the edit described by the diff:
# Insert two new rows at position 79 (pushing the existing rows 79.. down by 2),
# and populate the two new rows with new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before (old) row 79.
$ws.Rows.Item(79).Insert()
$ws.Rows.Item(79).Insert()

# --- New row 79 ---
$ws.Range("A79").Value = 9
$ws.Range("B79").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C79").Value = "Metropolitana"
$ws.Range("D79").Value = 44566
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100101
$ws.Range("H79").Value = "Berries"
$ws.Range("I79").Value = 100101001
$ws.Range("J79").Value = "Arándano (blue)"
$ws.Range("K79").Value = "Sin especificar"
$ws.Range("L79").Value = "Especial"
$ws.Range("M79").Value = 330
$ws.Range("N79").Value = 5000
$ws.Range("O79").Value = 5000
$ws.Range("P79").Value = 5000
$ws.Range("Q79").Value = '$/bandeja 2 kilos'
$ws.Range("R79").Value = "Provincia de Cardenal Caro"
$ws.Range("S79").Value = 2500
$ws.Range("T79").Value = 2

# --- New row 80 ---
$ws.Range("A80").Value = 9
$ws.Range("B80").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C80").Value = "Metropolitana"
$ws.Range("D80").Value = 44566
$ws.Range("E80").Value = 13
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100101
$ws.Range("H80").Value = "Berries"
$ws.Range("I80").Value = 100101001
$ws.Range("J80").Value = "Arándano (blue)"
$ws.Range("K80").Value = "Sin especificar"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 410
$ws.Range("N80").Value = 4000
$ws.Range("O80").Value = 4000
$ws.Range("P80").Value = 4000
$ws.Range("Q80").Value = '$/bandeja 2 kilos'
$ws.Range("R80").Value = "Provincia de Cardenal Caro"
$ws.Range("S80").Value = 2000
$ws.Range("T80").Value = 2
